$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column AV: "s_n_real_pred"
$ws.Range("AV1").Value = "s_n_real_pred"

# AV2: standalone formula (not part of the shared group below)
$ws.Range("AV2").Formula = "=(AQ2*X2-2*X2*(1-0.01*P2-2*0.01*AF2)/(-0.08/0.4*0.01*P2-(2*0.08/0.4+3)*0.01*AF2+1+0.08/0.4)+4*0.5*232000*(-0.4*0.01*AF2-0.08*0.01*P2)/((1-2*0.01*AF2)*0.4*0.08))/1000"

# AV3:AV18 share one formula (relative references fill down automatically)
$ws.Range("AV3:AV18").Formula = "=(AQ3*X3-2*X3*(1-0.01*P3-2*0.01*AF3)/(-0.08/0.4*0.01*P3-(2*0.08/0.4+3)*0.01*AF3+1+0.08/0.4)+4*0.5*232000*(-0.4*0.01*AF3-0.08*0.01*P3)/((1-2*0.01*AF3)*0.4*0.08))/1000"

$ws.Range("AV2").Select()
